{"js": "// Apply the \"Surat Bebas Pustaka\" edits:\n//  - Nomor number changed\n//  - Nama / NRP / Departemen values for the student changed\n\nconst body = context.document.body;\n\nconst replacements = [\n  {\n    find: \"Nomor : 1001/EBP ITS/6/2022\",\n    replace: \"Nomor : 57/EBP ITS/6/2022\"\n  },\n  {\n    find: \"         Nama           : Patrick\",\n    replace: \"         Nama           : Alexie Price MD\"\n  },\n  {\n    find: \"         NRP             : 05111840000098\",\n    replace: \"         NRP             : 04111340000410\"\n  },\n  {\n    find: \"         Departemen : Informatika\",\n    replace: \"         Departemen : expedite mission-critical infrastructures\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the \"Surat Bebas Pustaka\" edits:\n#  - Nomor number changed\n#  - Nama / NRP / Departemen values for the student changed\n\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphText($matchSubstring, $newText) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -like \"*$matchSubstring*\") {\n            $p.Range.Text = $newText\n            return $true\n        }\n    }\n    throw \"Paragraph not found containing: $matchSubstring\"\n}\n\nSet-ParagraphText \"Nomor : 1001/EBP ITS/6/2022\" \"Nomor : 57/EBP ITS/6/2022\" | Out-Null\nSet-ParagraphText \"Nama           : Patrick\" \"         Nama           : Alexie Price MD\" | Out-Null\nSet-ParagraphText \"NRP             : 05111840000098\" \"         NRP             : 04111340000410\" | Out-Null\nSet-ParagraphText \"Departemen : Informatika\" \"         Departemen : expedite mission-critical infrastructures\" | Out-Null\n"}
